$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'face/face082.png'
$ws.Cells.Item(2, 3).Value = 'formen'
$ws.Cells.Item(2, 4).Value = 'face'
$ws.Cells.Item(3, 2).Value = 'flower/flower083.png'
$ws.Cells.Item(3, 3).Value = 'stechen'
$ws.Cells.Item(3, 4).Value = 'flower'
$ws.Cells.Item(4, 2).Value = 'face/face068.png'
$ws.Cells.Item(4, 3).Value = 'tauschen'
$ws.Cells.Item(4, 4).Value = 'face'
$ws.Cells.Item(5, 2).Value = 'face/face112.png'
$ws.Cells.Item(5, 3).Value = 'pflegen'
$ws.Cells.Item(5, 4).Value = 'face'
$ws.Cells.Item(6, 2).Value = 'flower/flower086.png'
$ws.Cells.Item(6, 3).Value = 'spielen'
$ws.Cells.Item(6, 4).Value = 'flower'
$ws.Cells.Item(7, 2).Value = 'face/face072.png'
$ws.Cells.Item(7, 3).Value = 'ehren'
$ws.Cells.Item(7, 4).Value = 'face'
$ws.Cells.Item(8, 2).Value = 'flower/flower104.png'
$ws.Cells.Item(8, 3).Value = 'runden'
$ws.Cells.Item(8, 4).Value = 'flower'
$ws.Cells.Item(9, 2).Value = 'flower/flower092.png'
$ws.Cells.Item(9, 3).Value = 'hoffen'
$ws.Cells.Item(9, 4).Value = 'flower'
$ws.Cells.Item(10, 2).Value = 'flower/flower074.png'
$ws.Cells.Item(10, 3).Value = 'sieben'
$ws.Cells.Item(10, 4).Value = 'flower'
$ws.Cells.Item(11, 2).Value = 'flower/flower117.png'
$ws.Cells.Item(11, 3).Value = 'mieten'
$ws.Cells.Item(11, 4).Value = 'flower'
$ws.Cells.Item(12, 2).Value = 'face/face090.png'
$ws.Cells.Item(12, 3).Value = 'bleiben'
$ws.Cells.Item(12, 4).Value = 'face'
$ws.Cells.Item(13, 2).Value = 'flower/flower118.png'
$ws.Cells.Item(13, 3).Value = 'drehen'
$ws.Cells.Item(13, 4).Value = 'flower'
$ws.Cells.Item(14, 2).Value = 'face/face076.png'
$ws.Cells.Item(14, 3).Value = 'loben'
$ws.Cells.Item(14, 4).Value = 'face'
$ws.Cells.Item(15, 2).Value = 'flower/flower090.png'
$ws.Cells.Item(15, 3).Value = 'saufen'
$ws.Cells.Item(15, 4).Value = 'flower'
$ws.Cells.Item(16, 2).Value = 'flower/flower100.png'
$ws.Cells.Item(16, 3).Value = 'fesseln'
$ws.Cells.Item(16, 4).Value = 'flower'
$ws.Cells.Item(17, 2).Value = 'face/face100.png'
$ws.Cells.Item(17, 3).Value = 'lehnen'
$ws.Cells.Item(17, 4).Value = 'face'
$ws.Cells.Item(18, 2).Value = 'face/face103.png'
$ws.Cells.Item(18, 3).Value = 'bitten'
$ws.Cells.Item(18, 4).Value = 'face'
$ws.Cells.Item(19, 2).Value = 'face/face086.png'
$ws.Cells.Item(19, 3).Value = 'strahlen'
$ws.Cells.Item(19, 4).Value = 'face'
$ws.Cells.Item(20, 2).Value = 'face/face083.png'
$ws.Cells.Item(20, 3).Value = 'raten'
$ws.Cells.Item(20, 4).Value = 'face'
$ws.Cells.Item(21, 2).Value = 'flower/flower068.png'
$ws.Cells.Item(21, 3).Value = 'währen'
$ws.Cells.Item(21, 4).Value = 'flower'
$ws.Cells.Item(22, 2).Value = 'flower/flower065.png'
$ws.Cells.Item(22, 3).Value = 'fühlen'
$ws.Cells.Item(22, 4).Value = 'flower'
$ws.Cells.Item(23, 2).Value = 'face/face066.png'
$ws.Cells.Item(23, 3).Value = 'gelten'
$ws.Cells.Item(23, 4).Value = 'face'
$ws.Cells.Item(24, 2).Value = 'face/face095.png'
$ws.Cells.Item(24, 3).Value = 'hupen'
$ws.Cells.Item(24, 4).Value = 'face'
$ws.Cells.Item(25, 2).Value = 'flower/flower109.png'
$ws.Cells.Item(25, 3).Value = 'nehmen'
$ws.Cells.Item(25, 4).Value = 'flower'
$ws.Cells.Item(26, 2).Value = 'flower/flower064.png'
$ws.Cells.Item(26, 3).Value = 'schmecken'
$ws.Cells.Item(26, 4).Value = 'flower'
$ws.Cells.Item(27, 2).Value = 'face/face065.png'
$ws.Cells.Item(27, 3).Value = 'schenken'
$ws.Cells.Item(27, 4).Value = 'face'
$ws.Cells.Item(28, 2).Value = 'face/face119.png'
$ws.Cells.Item(28, 3).Value = 'husten'
$ws.Cells.Item(28, 4).Value = 'face'
$ws.Cells.Item(29, 2).Value = 'flower/flower108.png'
$ws.Cells.Item(29, 3).Value = 'scheitern'
$ws.Cells.Item(29, 4).Value = 'flower'
$ws.Cells.Item(30, 2).Value = 'face/face085.png'
$ws.Cells.Item(30, 3).Value = 'wiegen'
$ws.Cells.Item(30, 4).Value = 'face'
$ws.Cells.Item(31, 2).Value = 'face/face093.png'
$ws.Cells.Item(31, 3).Value = 'kehren'
$ws.Cells.Item(31, 4).Value = 'face'
$ws.Cells.Item(32, 2).Value = 'flower/flower098.png'
$ws.Cells.Item(32, 3).Value = 'hauen'
$ws.Cells.Item(32, 4).Value = 'flower'
$ws.Cells.Item(33, 2).Value = 'flower/flower077.png'
$ws.Cells.Item(33, 3).Value = 'gründen'
$ws.Cells.Item(33, 4).Value = 'flower'
